# Update "想去人数" (F column) values across sheets to match the
# newly generated output snapshot (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(9, 6).Value  = 9943
$ws.Cells.Item(11, 6).Value = 2968
$ws.Cells.Item(14, 6).Value = 2724
$ws.Cells.Item(17, 6).Value = 2118
$ws.Cells.Item(20, 6).Value = 375
$ws.Cells.Item(25, 6).Value = 184
$ws.Cells.Item(27, 6).Value = 1303
$ws.Cells.Item(32, 6).Value = 2039
$ws.Cells.Item(33, 6).Value = 2889
$ws.Cells.Item(34, 6).Value = 6
$ws.Cells.Item(37, 6).Value = 369
$ws.Cells.Item(40, 6).Value = 70

# --- 演出 (Performances) sheet ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(15, 6).Value = 166

# --- 本地生活 (Local Life) sheet ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 119

# --- 全部类型 (All Types) sheet ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(4, 6).Value  = 119
$ws.Cells.Item(13, 6).Value = 9943
$ws.Cells.Item(16, 6).Value = 2969
$ws.Cells.Item(19, 6).Value = 2724
$ws.Cells.Item(21, 6).Value = 2118
$ws.Cells.Item(24, 6).Value = 375
$ws.Cells.Item(28, 6).Value = 184
$ws.Cells.Item(30, 6).Value = 1303
$ws.Cells.Item(34, 6).Value = 2040
$ws.Cells.Item(36, 6).Value = 2889
$ws.Cells.Item(39, 6).Value = 369
$ws.Cells.Item(45, 6).Value = 70
$ws.Cells.Item(49, 6).Value = 166
